$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 643
$ws1.Range("F4").Value = 1446
$ws1.Range("F5").Value = 682

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 643
$ws4.Range("F4").Value = 1446
$ws4.Range("F6").Value = 682
